# The post originally at row 169 ("「この世界で君も何かを持っている。だから立ち
# 上がれ！」" / Ghassan Kanafani quote) was removed from the posts sheet.
# Deleting the entire row shifts every following row up by one, which is
# exactly what the commit's row-by-row renumbering (169->dropped, 170->169,
# 171->170, ... 245->244) reflects, and updates the sheet's used range
# from A1:C245 to A1:C244 automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(169).Delete()
